$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Rename the sheet from "其他有價證券" to "具有相當價值之財產"
$ws.Name = "具有相當價值之財產"

# Remove the two old data/sample rows (rows 3 and 4); only a header row (1)
# and a single data row (2) remain in the new layout.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Extend formatting (style) for the new columns H:L on row 1 and row 2 by
# copying the existing formatted cells before overwriting their values.
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,8))
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,9))
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,10))
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,11))
$ws.Cells.Item(1,3).Copy($ws.Cells.Item(1,12))

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,8))
$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,9))
$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,10))
$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,11))
$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,12))

# Header row (row 1)
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "quantity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"
$ws.Cells.Item(1,6).Value = "property_category"
$ws.Cells.Item(1,7).Value = "category"
$ws.Cells.Item(1,8).Value = "date"
$ws.Cells.Item(1,9).Value = "legislator_name"
$ws.Cells.Item(1,10).Value = "legislator_id"
$ws.Cells.Item(1,11).Value = "source_file"
$ws.Cells.Item(1,12).Value = "index"

# Data row (row 2)
$ws.Cells.Item(2,1).Value = 128
$ws.Cells.Item(2,2).Value = "朱銘雕刻"
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = "陳明文"
$ws.Cells.Item(2,5).Value = 500000
$ws.Cells.Item(2,6).Value = "otherbonds"
$ws.Cells.Item(2,7).Value = "normal"
$ws.Cells.Item(2,8).Value = "2012-02-13"
$ws.Cells.Item(2,9).Value = "陳明文"
$ws.Cells.Item(2,10).Value = 828
$ws.Cells.Item(2,11).Value = "tmpf4561"
$ws.Cells.Item(2,12).Value = 128
